# UML Merge & Modify Actor Description
#
# The "Actor descriptions" table (rows 12-17, columns C/D) is collapsed
# from 6 rows down to 4 rows:
#   - "고객" + "회원"(old) are merged into a single "회원" actor whose
#     description now covers sign-up, login and the clothing actions.
#   - The old "고객(판매자)" row becomes a new "프로그램 종료" actor that
#     describes automatic program termination on logout.
#   - "고객(구매자)" is dropped; "택배사" and "이메일" shift up to take
#     its place, becoming the last two rows of the table.
#
# Because the table shrinks by two rows, the former last row (row 17,
# which carries the "closing border" style) has to become row 15, so we
# copy that formatting onto row 15 before wiping out rows 16-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Pull the bottom-border formatting of the old last row (17) onto the
#    row that will become the new last row (15) of the actor table.
$ws.Range("C17:E17").Copy()
$ws.Range("C15:E15").PasteSpecial(-4122)

# 2) Rewrite the actor/description pairs in their new, merged form.
$ws.Range("C12").Value = "회원"
$ws.Range("D12").Value = "회원은 회원가입을 할 수 있고, 로그인을 통해 시스템에 대한 사용 권한을 얻을 수 있습니다. 또한 로그인 하여 의류를 등록/ 판매/ 조회 할 수 있습니다."

$ws.Range("C13").Value = "프로그램 종료"
$ws.Range("D13").Value = "회원이 로그아웃 시 자동으로 프로그램을 종료합니다."

$ws.Range("C14").Value = "택배사"
$ws.Range("D14").Value = "..."

$ws.Range("C15").Value = "이메일"
$ws.Range("D15").Value = "..."

# 3) The table is now only 4 rows (12-15); drop the two now-unused rows
#    entirely (full clear, not just contents) so they vanish from the
#    sheet rather than lingering as blank formatted rows.
$ws.Range("C16:E17").Clear()
